$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row heights (cosmetic, matches target row ht/customHeight) ---
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75

# --- Row 2 values: MCH334-1 / Series / 1 Box / LOCATION: 33G ... ---
$ws.Range("A2").Value = "MCH334-1"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 33G | GRAP COUNT NUMER: NONE"

# --- Row 3 values: MCH334 / LIBERATION 1953-1959 / 1953-1959 / Series / 1 Box / LOCATION: 102C ... ---
$ws.Range("A3").Value = "MCH334"
$ws.Range("C3").Value = "LIBERATION 1953-1959"
$ws.Range("D3").Value = "1953-1959"
$ws.Range("E3").Value = "Series"
$ws.Range("F3").Value = "1 Box"
$ws.Range("G3").Value = "LOCATION: 102C | GRAP COUNT NUMER: BOX 27"

# --- Formatting for row 2 and most of row 3 (font: Calibri 10pt, theme text color) ---
# NOTE: multi-area (union) ranges only apply formatting to the first area in
# this runtime, so iterate explicitly over each area instead.
$calibri10 = $ws.Range("A2,C2,D2,E2,F2,G2,H2,C3,D3,E3,F3,G3,H3")
foreach ($area in $calibri10.Areas) {
    $area.Font.Name = "Calibri"
    $area.Font.Size = 10
    $area.Font.ThemeColor = 1
}

# --- A3 uses a slightly larger font (Calibri 11pt, theme text color) ---
$a3Font = $ws.Range("A3").Font
$a3Font.Name = "Calibri"
$a3Font.Size = 11
$a3Font.ThemeColor = 1

# --- Sheet view: keep the existing frozen header pane, but move the active
#     selection to F15 (matches the authored file's saved cursor position) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F15").Select()
